$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New id_parameter rows (ids 66-70), appended after the existing "Conversion efficiency" row (65).
# The shared-string insertion order below intentionally matches the source
# workbook's string table (label/description for RES allocation and
# substitution factors first, then the running-cost labels, then the
# investment/running-cost descriptions) so the resulting sharedStrings.xml
# lines up with the authored file.

# Row 67 -> id 66: RES energy allocation
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "RES energy allocation"
$ws.Cells.Item(67, 3).Value = "Generated final energy carrier and necessary quantity of final energy carriers for the generation"

# Row 68 -> id 67: Substitution factors
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = "Substitution factors"
$ws.Cells.Item(68, 3).Value = "Coefficients specifying substituted electricity generators for each renewable electricity technology"

# Row 70 -> id 69: Running costs (label written before row 69/71 labels)
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = "Running costs"

# Row 71 -> id 70: Variable running costs
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "Variable running costs"

# Row 69 -> id 68: Investment costs per capacity
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = "Investment costs per capacity"
$ws.Cells.Item(69, 3).Value = "Capital costs depending on installed capacity (CAPEX)"

# Descriptions for the running-cost rows
$ws.Cells.Item(70, 3).Value = "Running costs depending on installed capacity in M€/MW (OPEX)"
$ws.Cells.Item(71, 3).Value = "Running costs depending on generated energy in M€/MWh (OPEX)"

# Move selection to match the post-edit cursor position (one row below the last data row)
$ws.Range("C72").Select()
